# Updated cryptos list (Price and Volume(1h) columns) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" values are plain decimal-looking strings (e.g. "89.50"); Excel
# auto-converts these to numbers on assignment (losing formatting / precision), so
# force a text number-format before writing, then restore the default cell style so
# no stray formatting is left behind.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.163.66'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.46%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.871.71'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.89%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5028'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3747'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07155'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8903'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.14%  '
$ws.Range('E11').Value = '  -0.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.886.38'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07574'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.329'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.50'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008513'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.15'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.22%  '
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.218.20'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.086'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.131.35'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.63'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.503'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.16'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.68%  '
$ws.Range('E26').Value = '  -2.33%  '
$ws.Range('E27').Value = '  -2.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.090'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.01'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.770'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.704'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08976'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05148'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.095'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7478'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.164'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.565'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02034'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.044'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.60%  '
$ws.Range('E40').Value = '  -1.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5364'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.640'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '115.32'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.476'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1480'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4650'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.001'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('E48').Value = '  -4.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.574'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '64.85'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.82'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.45%  '
